$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.54
$ws.Range("G2").Value = 1.61
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 4.5
$ws.Range("K2").Value = 5.1
$ws.Range("L2").Value = 1.36
$ws.Range("N2").Value = 4.3
$ws.Range("Q2").Value = 1.8
$ws.Range("S2").Value = 3.1
$ws.Range("T2").Value = 1.89
$ws.Range("U2").Value = 1.94
$ws.Range("W2").Value = 2.62
$ws.Range("Y2").Value = 55
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 14
$ws.Range("AD2").Value = 990
$ws.Range("AG2").Value = 20
$ws.Range("AH2").Value = 65
$ws.Range("AI2").Value = 330
$ws.Range("AJ2").Value = 40
$ws.Range("AL2").Value = 160
$ws.Range("AN2").Value = 10
$ws.Range("F3").Value = 7.4
$ws.Range("G3").Value = 8.4
$ws.Range("H3").Value = 1.42
$ws.Range("I3").Value = 1.45
$ws.Range("K3").Value = 6
$ws.Range("P3").Value = 2.48
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.58
$ws.Range("S3").Value = 2.52
$ws.Range("T3").Value = 1.84
$ws.Range("X3").Value = 27
$ws.Range("AB3").Value = 85
$ws.Range("AH3").Value = 46
$ws.Range("AI3").Value = 55
$ws.Range("F4").Value = 1.44
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = 6.8
$ws.Range("I4").Value = 28
$ws.Range("J4").Value = 4.6
$ws.Range("K4").Value = 6.2
$ws.Range("L4").Value = 1.32
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 4.9
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 2.68
$ws.Range("V4").Value = 1.11
$ws.Range("W4").Value = 2.96
$ws.Range("AC4").Value = 990
$ws.Range("AJ4").Value = 900
$ws.Range("F5").Value = 2.18
$ws.Range("H5").Value = 3.55
$ws.Range("L5").Value = 1.46
$ws.Range("N5").Value = 3.2
$ws.Range("O5").Value = 1.38
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.16
$ws.Range("S5").Value = 3.95
$ws.Range("T5").Value = 1.85
$ws.Range("W5").Value = 1.73
$ws.Range("X5").Value = 12.5
$ws.Range("Y5").Value = 13.5
$ws.Range("Z5").Value = 80
$ws.Range("AB5").Value = 10
$ws.Range("AC5").Value = 8.6
$ws.Range("AD5").Value = 18.5
$ws.Range("AF5").Value = 14.5
$ws.Range("AH5").Value = 21
$ws.Range("AJ5").Value = 32
$ws.Range("AK5").Value = 30
$ws.Range("AO5").Value = 70
$ws.Range("F6").Value = 1.79
$ws.Range("G6").Value = 1.89
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 5.9
$ws.Range("J6").Value = 3.45
$ws.Range("L6").Value = 1.48
$ws.Range("N6").Value = 3.25
$ws.Range("O6").Value = 1.38
$ws.Range("P6").Value = 1.77
$ws.Range("Q6").Value = 2.14
$ws.Range("R6").Value = 1.28
$ws.Range("S6").Value = 3.95
$ws.Range("T6").Value = 1.92
$ws.Range("U6").Value = 1.89
$ws.Range("W6").Value = 2.12
$ws.Range("Y6").Value = 30
$ws.Range("AB6").Value = 8
$ws.Range("AD6").Value = 46
$ws.Range("AF6").Value = 11.5
$ws.Range("AN6").Value = 42
